$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: -4,-8)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: -5,-9)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: 0,8)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: -1,-1)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: -2,-8)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: 6,5)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: -6,0)"

$ws.Range("A3").Value = "cost: 434.7749194859469"
$ws.Range("A4").Value = "time: 81.95498389718938"
